# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with refreshed quote data, matching the "Updated cryptos list ..." commit.
# Price values are stored as text in the sheet (e.g. "27.091.61"), so a
# leading apostrophe is used to force text entry and stop Excel from
# auto-converting them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.091.61"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "'1.823.79"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'311.56"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "'0.4671"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").Value = "'0.3643"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").Value = "'0.07299"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").Value = "'0.8697"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "'20.18"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'0.07616"
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("D13").Value = "'1.846.47"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'93.03"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'5.343"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "'6.475"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "'0.000008652"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D20").Value = "'27.284.74"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "'14.50"
$ws.Range("E21").Value = "  -2.37%  "
$ws.Range("D22").Value = "'5.194"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("D24").Value = "'2.072.28"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "'151.73"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'1.857"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").Value = "'18.27"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").Value = "'2.112"
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").Value = "'5.098"
$ws.Range("E29").Value = "  -3.47%  "
$ws.Range("D30").Value = "'115.97"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "'0.08931"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "'0.7344"
$ws.Range("E33").Value = "  -3.62%  "
$ws.Range("D34").Value = "'4.457"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("D35").Value = "'1.143"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").Value = "'1.009"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").Value = "'2.547"
$ws.Range("E37").Value = "  +6.86%  "
$ws.Range("D38").Value = "'0.05270"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").Value = "'1.071"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("D40").Value = "'0.01920"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "'2.935"
$ws.Range("E41").Value = "  -2.25%  "
$ws.Range("D42").Value = "'7.123"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("D43").Value = "'0.5228"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").Value = "'8.274"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("D46").Value = "'0.4875"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "'1.009"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "'103.78"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").Value = "'10.13"
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("D50").Value = "'1.639"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("D51").Value = "'0.06249"
